# "write memo-box css in the favorite section"
# Slide with sldId=261 (the 6th slide in the deck) has an ellipse shape
# (id=9, name "楕円 8", creationId {05A3C300-86CA-4E40-9F97-058B51978EC2})
# used as the background of the favorite/memo section. Its outline was
# "No Line"; give it a solid outline matching the shape's own fill color
# (2B8DF9) so it reads as a bordered memo box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)

$shp.Line.Visible = $true
$shp.Line.ForeColor.RGB = 16354603   # 0x2B8DF9 (R=43,G=141,B=249) packed as R+G*256+B*65536
